# Finishing the data mapper: add the BCA master-data row to the hrd_bank sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hrd_bank")

$ws.Range("A3").Value = "BCA"
$ws.Range("B3").Value = "Bank Central Asia"

# Match the existing data rows' explicit row height so row 3 is formatted like rows 1-2.
$ws.Rows.Item(3).RowHeight = 16.5

# Leave the cursor where the user would land after typing the last cell and pressing Enter/Tab.
$ws.Range("B4").Select() | Out-Null
